$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 123..127 (match data in columns F..V) got cyclically re-ordered.
#    Column A (Indice) stays a strict sequence and columns A..E (Indice,
#    pais, torneio, temporada, data_partida) are unchanged for each row;
#    only the match info (home/away teams, goals, odds, timestamps, url)
#    moved between rows.  Capture all five rows first, then write them
#    back in their new positions so the copy cannot clobber itself.
# ---------------------------------------------------------------------------
$row123 = $ws.Range("F123:V123").Value2
$row124 = $ws.Range("F124:V124").Value2
$row125 = $ws.Range("F125:V125").Value2
$row126 = $ws.Range("F126:V126").Value2
$row127 = $ws.Range("F127:V127").Value2

$ws.Range("F123:V123").Value2 = $row125
$ws.Range("F124:V124").Value2 = $row123
$ws.Range("F125:V125").Value2 = $row124
$ws.Range("F126:V126").Value2 = $row127
$ws.Range("F127:V127").Value2 = $row126

# ---------------------------------------------------------------------------
# 2) A brand-new match (Fulham vs Wolves) was appended as row 131, pushing
#    the sheet dimension from A1:V130 to A1:V131.
# ---------------------------------------------------------------------------
$ws.Range("A130:V130").Copy()
$ws.Range("A131:V131").PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(131, 1).Value = 130
$ws.Cells.Item(131, 2).Value = "england"
$ws.Cells.Item(131, 3).Value = "premier-league"
$ws.Cells.Item(131, 4).Value = "2023-2024"
$ws.Cells.Item(131, 5).Value = 45257.875
$ws.Cells.Item(131, 6).Value = "Fulham"
$ws.Cells.Item(131, 7).Value = 3
$ws.Cells.Item(131, 8).Value = "Wolves"
$ws.Cells.Item(131, 9).Value = 2
$ws.Cells.Item(131, 10).Value = 2.04
$ws.Cells.Item(131, 11).Value = "05/11/2023 11:03"
$ws.Cells.Item(131, 12).Value = 2.43
$ws.Cells.Item(131, 13).Value = "27/11/2023 20:50"
$ws.Cells.Item(131, 14).Value = 3.48
$ws.Cells.Item(131, 15).Value = "05/11/2023 11:03"
$ws.Cells.Item(131, 16).Value = 3.3
$ws.Cells.Item(131, 17).Value = "27/11/2023 20:57"
$ws.Cells.Item(131, 18).Value = 3.8
$ws.Cells.Item(131, 19).Value = "05/11/2023 11:03"
$ws.Cells.Item(131, 20).Value = 3.19
$ws.Cells.Item(131, 21).Value = "27/11/2023 20:57"
$ws.Cells.Item(131, 22).Value = "https://www.betexplorer.com/football/england/premier-league/fulham-wolves/x0jTpuq2/"
